# Update code and data
# Adds two small "pyramid by age-group" summary tables (A31:D36, G31:I36,
# M31:O36) below the existing detailed age-pyramid table on the single
# worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlCenter = -4108

# ---------------------------------------------------------------------
# 1. Header row (row 31) for the three new mini tables.
#    Style = same header look used in row 1 (bold font, fill, border).
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("A31,B31,C31,D31,G31,H31,I31,M31,N31,O31").PasteSpecial($xlPasteFormats)

$ws.Range("A31").Value = "Grupos de Edades"
$ws.Range("B31").Value = "Defunciones"
$ws.Range("D31").Value = "CFR Global"

# Age-group labels (also used in the blocks below) -- written here so the
# shared-string table picks them up before "% muertes por Grupos".
$ws.Range("A32").Value = "< 20 años"
$ws.Range("A33").Value = "20-39 años"
$ws.Range("A34").Value = "40-59 años"
$ws.Range("A35").Value = "60-79 años"
$ws.Range("A36").Value = "> 80 años"

$ws.Range("C31").Value = "% muertes por Grupos"

$ws.Range("G31").Value = "Grupos de Edades"
$ws.Range("H31").Value = $ws.Range("F1").Value
$ws.Range("I31").Value = $ws.Range("I1").Value

$ws.Range("M31").Value = "Grupos de Edades"
$ws.Range("N31").Value = $ws.Range("D1").Value
$ws.Range("O31").Value = $ws.Range("E1").Value

# ---------------------------------------------------------------------
# 2. Body rows (32-36) base formatting: thin box border + centered text,
#    same look used throughout the sheet (copy border from an existing
#    bordered cell, then center).
# ---------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("A32:A36,B32:B36,G32:G36,H32:H36,M32:M36,N32:N36").PasteSpecial($xlPasteFormats)
$ws.Range("A32:A36,B32:B36,G32:G36,H32:H36,M32:M36,N32:N36").HorizontalAlignment = $xlCenter

$ws.Range("B2").Copy()
$ws.Range("C32:C36,I32:I36,O32:O36").PasteSpecial($xlPasteFormats)
$ws.Range("C32:C36,I32:I36,O32:O36").HorizontalAlignment = $xlCenter
$ws.Range("C32:C36,I32:I36,O32:O36").NumberFormat = "0%"

$ws.Range("B2").Copy()
$ws.Range("D32:D36").PasteSpecial($xlPasteFormats)
$ws.Range("D32:D36").HorizontalAlignment = $xlCenter
$ws.Range("D32:D36").NumberFormat = "0.00%"

# ---------------------------------------------------------------------
# 3. Block A (A:D) -- Defunciones (deaths) per age-group, with CFR-style
#    percentages computed from the grand totals in row 24.
# ---------------------------------------------------------------------
$ws.Range("B32").Value = 27
$ws.Range("B33").Value = 186
$ws.Range("B34").Value = 794
$ws.Range("B35").Value = 1970
$ws.Range("B36").Value = 1045

$ws.Range("C32").Formula = "=B32/G24"
$ws.Range("C33").Formula = "=B33/G24"
$ws.Range("C34").Formula = "=B34/G24"
$ws.Range("C35").Formula = "=B35/G24"
$ws.Range("C36").Formula = "=B36/G24"

$ws.Range("D32").Formula = "=B32/F24"
$ws.Range("D33").Formula = "=B33/F24"
$ws.Range("D34").Formula = "=B34/F24"
$ws.Range("D35").Formula = "=B35/F24"
$ws.Range("D36").Formula = "=B36/F24"

# ---------------------------------------------------------------------
# 4. Block B (G:I) -- Casos Covid (cases) per age-group.
# ---------------------------------------------------------------------
$ws.Range("G32").Value = "< 20 años"
$ws.Range("G33").Value = "20-39 años"
$ws.Range("G34").Value = "40-59 años"
$ws.Range("G35").Value = "60-79 años"
$ws.Range("G36").Value = "> 80 años"

$ws.Range("H32").Value = 41893
$ws.Range("H33").Value = 99101
$ws.Range("H34").Value = 72386
$ws.Range("H35").Value = 27965
$ws.Range("H36").Value = 5445

$ws.Range("I32").Value = 0.17
$ws.Range("I33").Value = 0.4
$ws.Range("I34").Value = 0.29
$ws.Range("I35").Value = 0.11
$ws.Range("I36").Value = 0.02

# ---------------------------------------------------------------------
# 5. Block C (M:O) -- Población (population) per age-group.
# ---------------------------------------------------------------------
$ws.Range("M32").Value = "< 20 años"
$ws.Range("M33").Value = "20-39 años"
$ws.Range("M34").Value = "40-59 años"
$ws.Range("M35").Value = "60-79 años"
$ws.Range("M36").Value = "> 80 años"

$ws.Range("N32").Value = 1491818
$ws.Range("N33").Value = 1280979
$ws.Range("N34").Value = 967156
$ws.Range("N35").Value = 422595
$ws.Range("N36").Value = 83892

$ws.Range("O32").Value = 0.35
$ws.Range("O33").Value = 0.3
$ws.Range("O34").Value = 0.23
$ws.Range("O35").Value = 0.1
$ws.Range("O36").Value = 0.02

# ---------------------------------------------------------------------
# 6. Column widths for the newly-used columns.
# ---------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 14.997679337211277
$ws.Columns("C:C").ColumnWidth = 18.828408946161684
$ws.Columns("G:G").ColumnWidth = 14.997679337211277
$ws.Columns("I:I").ColumnWidth = 16.49716982634171
$ws.Columns("M:M").ColumnWidth = 14.997679337211277
$ws.Columns("O:O").ColumnWidth = 16.49716982634171

# ---------------------------------------------------------------------
# 7. Selection, matching the author's last recorded cursor position.
# ---------------------------------------------------------------------
$ws.Range("O32:O34").Select()
